$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '34.414.41'
$ws.Cells.Item(2, 5).Value = '  -0.66%  '
$ws.Cells.Item(3, 4).Value = '1.804.47'
$ws.Cells.Item(3, 5).Value = '  -0.48%  '
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '227.64'
$ws.Cells.Item(5, 5).Value = '  -0.29%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.604'
$ws.Cells.Item(6, 5).Value = '  +6.67%  '
$ws.Cells.Item(7, 5).Value = '  +0.11%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '36.29'
$ws.Cells.Item(8, 5).Value = '  +4.11%  '
$ws.Cells.Item(9, 5).Value = '  -0.58%  '
$ws.Cells.Item(10, 5).Value = '  +0.37%  '
$ws.Cells.Item(11, 5).Value = '  +1.33%  '
$ws.Cells.Item(12, 4).Value = '2.066.18'
$ws.Cells.Item(12, 5).Value = '  -0.38%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '11.29'
$ws.Cells.Item(13, 5).Value = '  -0.26%  '
$ws.Cells.Item(14, 4).Value = '1.811.06'
$ws.Cells.Item(14, 5).Value = '  +1.08%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.647'
$ws.Cells.Item(15, 5).Value = '  +0.60%  '
$ws.Cells.Item(16, 5).Value = '  +2.98%  '
$ws.Cells.Item(17, 4).Value = '34.420.40'
$ws.Cells.Item(17, 5).Value = '  -0.70%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '69.79'
$ws.Cells.Item(18, 5).Value = '  +1.19%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '245.88'
$ws.Cells.Item(19, 5).Value = '  -0.68%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0788'
$ws.Cells.Item(20, 5).Value = '  -1.60%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '11.49'
$ws.Cells.Item(21, 5).Value = '  -0.24%  '
$ws.Cells.Item(22, 5).Value = '  +0.06%  '
$ws.Cells.Item(23, 5).Value = '  -0.06%  '
$ws.Cells.Item(24, 5).Value = '  +6.65%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '171.27'
$ws.Cells.Item(25, 5).Value = '  -0.30%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.04'
$ws.Cells.Item(26, 5).Value = '  +8.03%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '17.54'
$ws.Cells.Item(27, 5).Value = '  +4.80%  '
$ws.Cells.Item(28, 5).Value = '  +4.18%  '
$ws.Cells.Item(29, 5).Value = '  +0.02%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '4.00'
$ws.Cells.Item(30, 5).Value = '  +0.57%  '
$ws.Cells.Item(31, 5).Value = '  +0.03%  '
$ws.Cells.Item(32, 5).Value = '  -0.66%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0524'
$ws.Cells.Item(33, 5).Value = '  -1.52%  '
$ws.Cells.Item(34, 5).Value = '  -2.97%  '
$ws.Cells.Item(35, 4).Value = '1.382.93'
$ws.Cells.Item(35, 5).Value = '  -2.71%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.661'
$ws.Cells.Item(36, 5).Value = '  -2.74%  '
$ws.Cells.Item(37, 5).Value = '  -0.84%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.37'
$ws.Cells.Item(38, 5).Value = '  -10.43%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.83'
$ws.Cells.Item(40, 5).Value = '  -0.94%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '82.18'
$ws.Cells.Item(41, 5).Value = '  -3.70%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.953'
$ws.Cells.Item(42, 5).Value = '  -0.23%  '
$ws.Cells.Item(43, 5).Value = '  +0.58%  '
$ws.Cells.Item(44, 5).Value = '  +7.52%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.54'
$ws.Cells.Item(45, 5).Value = '  -2.31%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '5.97'
$ws.Cells.Item(46, 5).Value = '  -2.21%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.0500'
$ws.Cells.Item(47, 5).Value = '  -4.02%  '
$ws.Cells.Item(48, 4).Value = '1.967.28'
$ws.Cells.Item(48, 5).Value = '  -0.46%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.00'
$ws.Cells.Item(49, 5).Value = '  +0.03%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '102.61'
$ws.Cells.Item(50, 5).Value = '  -2.55%  '
$ws.Cells.Item(51, 5).Value = '  -2.43%  '
